$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 261, shifting existing rows 261-272 down to 262-273
$ws.Rows.Item(261).Insert()

# Copy the date number format used by the other rows in column D onto the new cell
$ws.Cells.Item(261, 4).NumberFormat = $ws.Cells.Item(262, 4).NumberFormat

# Fill in the new row 261 with the new weekly record
$ws.Cells.Item(261, 1).Value = 7
$ws.Cells.Item(261, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(261, 3).Value = "Ñuble"
$ws.Cells.Item(261, 4).Value = 44753
$ws.Cells.Item(261, 5).Value = 16
$ws.Cells.Item(261, 6).Value = 100112023
$ws.Cells.Item(261, 7).Value = "Brócoli"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 200
$ws.Cells.Item(261, 11).Value = 800
$ws.Cells.Item(261, 12).Value = 900
$ws.Cells.Item(261, 13).Value = 850
$ws.Cells.Item(261, 14).Value = "`$/unidad"
$ws.Cells.Item(261, 15).Value = "Región del Maule"
$ws.Cells.Item(261, 16).Value = 850
$ws.Cells.Item(261, 17).Value = 1
$ws.Cells.Item(261, 18).Value = "Hortaliza"
